$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 0.001
$ws.Range("K3").Value = 684
$ws.Range("L3").Value = 0.00684
